{"js": "// Update the 25 \"three-digit \u00d7 one-digit\" answer cells in the practice\n// table to a new set of problems/answers (author regenerated the sheet).\n// Each old value is unique in the document, so a simple search/replace\n// per pair is sufficient and keeps every other run/paragraph untouched.\nconst replacements = [\n  [\"758\u00d79=6822\", \"212\u00d78=1696\"],\n  [\"759\u00d79=6831\", \"424\u00d79=3816\"],\n  [\"222\u00d75=1110\", \"616\u00d73=1848\"],\n  [\"609\u00d75=3045\", \"926\u00d73=2778\"],\n  [\"476\u00d76=2856\", \"207\u00d74=828\"],\n  [\"360\u00d76=2160\", \"582\u00d79=5238\"],\n  [\"832\u00d76=4992\", \"935\u00d78=7480\"],\n  [\"776\u00d79=6984\", \"817\u00d75=4085\"],\n  [\"177\u00d76=1062\", \"341\u00d72=682\"],\n  [\"600\u00d76=3600\", \"967\u00d74=3868\"],\n  [\"526\u00d72=1052\", \"153\u00d72=306\"],\n  [\"367\u00d79=3303\", \"318\u00d73=954\"],\n  [\"903\u00d76=5418\", \"390\u00d74=1560\"],\n  [\"562\u00d78=4496\", \"418\u00d78=3344\"],\n  [\"236\u00d73=708\", \"422\u00d78=3376\"],\n  [\"264\u00d78=2112\", \"429\u00d77=3003\"],\n  [\"681\u00d79=6129\", \"186\u00d72=372\"],\n  [\"149\u00d78=1192\", \"786\u00d74=3144\"],\n  [\"497\u00d72=994\", \"189\u00d76=1134\"],\n  [\"840\u00d78=6720\", \"333\u00d76=1998\"],\n  [\"199\u00d77=1393\", \"726\u00d78=5808\"],\n  [\"359\u00d77=2513\", \"940\u00d75=4700\"],\n  [\"305\u00d72=610\", \"871\u00d75=4355\"],\n  [\"863\u00d79=7767\", \"919\u00d74=3676\"],\n  [\"486\u00d72=972\", \"457\u00d72=914\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# Update the 25 \"three-digit x one-digit\" answer cells in the practice\n# table to a new set of problems/answers (author regenerated the sheet).\n# Each old value is unique in the document, so a simple Find/Replace\n# per pair is sufficient and keeps every other run/paragraph untouched.\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @{ Old = \"758\u00d79=6822\"; New = \"212\u00d78=1696\" },\n    @{ Old = \"759\u00d79=6831\"; New = \"424\u00d79=3816\" },\n    @{ Old = \"222\u00d75=1110\"; New = \"616\u00d73=1848\" },\n    @{ Old = \"609\u00d75=3045\"; New = \"926\u00d73=2778\" },\n    @{ Old = \"476\u00d76=2856\"; New = \"207\u00d74=828\" },\n    @{ Old = \"360\u00d76=2160\"; New = \"582\u00d79=5238\" },\n    @{ Old = \"832\u00d76=4992\"; New = \"935\u00d78=7480\" },\n    @{ Old = \"776\u00d79=6984\"; New = \"817\u00d75=4085\" },\n    @{ Old = \"177\u00d76=1062\"; New = \"341\u00d72=682\" },\n    @{ Old = \"600\u00d76=3600\"; New = \"967\u00d74=3868\" },\n    @{ Old = \"526\u00d72=1052\"; New = \"153\u00d72=306\" },\n    @{ Old = \"367\u00d79=3303\"; New = \"318\u00d73=954\" },\n    @{ Old = \"903\u00d76=5418\"; New = \"390\u00d74=1560\" },\n    @{ Old = \"562\u00d78=4496\"; New = \"418\u00d78=3344\" },\n    @{ Old = \"236\u00d73=708\";  New = \"422\u00d78=3376\" },\n    @{ Old = \"264\u00d78=2112\"; New = \"429\u00d77=3003\" },\n    @{ Old = \"681\u00d79=6129\"; New = \"186\u00d72=372\" },\n    @{ Old = \"149\u00d78=1192\"; New = \"786\u00d74=3144\" },\n    @{ Old = \"497\u00d72=994\";  New = \"189\u00d76=1134\" },\n    @{ Old = \"840\u00d78=6720\"; New = \"333\u00d76=1998\" },\n    @{ Old = \"199\u00d77=1393\"; New = \"726\u00d78=5808\" },\n    @{ Old = \"359\u00d77=2513\"; New = \"940\u00d75=4700\" },\n    @{ Old = \"305\u00d72=610\";  New = \"871\u00d75=4355\" },\n    @{ Old = \"863\u00d79=7767\"; New = \"919\u00d74=3676\" },\n    @{ Old = \"486\u00d72=972\";  New = \"457\u00d72=914\" }\n)\n\nforeach ($pair in $replacements) {\n    $range = $d.Content\n    $range.Find.Execute($pair.Old, $false, $false, $false, $false, $false, $true, 1, $false, $pair.New, 2) | Out-Null\n}\n"}
